# Update extreme case analysis
# 1. Sort the main data table on Sheet1 by Num of outages (B) descending,
#    then by Main Cause (C) ascending.
# 2. Add an "Average" summary row under the data on Sheet1.
# 3. Build two new sheets, "Trees" and "Lightning", containing the
#    (weather-complete) rows for each Main Cause, sorted by outages
#    descending, each with its own Average summary row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Sort the data block A1:N40 (header in row 1) by column B desc,
#    column C asc (tie-breaker).
# ---------------------------------------------------------------------
$dataRange = $ws.Range("A1:N40")
$key1 = $ws.Range("B1")
$key2 = $ws.Range("C1")
$dataRange.Sort($key1, 2, $key2, $null, 1, $null, $null, $null, $true, $null, $null, 1)

# ---------------------------------------------------------------------
# 2. Create the "Trees" and "Lightning" sheets, positioned after Sheet1.
# ---------------------------------------------------------------------
$treesSheet = $wb.Worksheets.Add($null, $ws)
$treesSheet.Name = "Trees"

$lightningSheet = $wb.Worksheets.Add($null, $treesSheet)
$lightningSheet.Name = "Lightning"

# ---------------------------------------------------------------------
# Walk the sorted Sheet1 rows (2..40) once, copying each "complete"
# weather row (column F populated) into the Trees or Lightning sheet
# depending on its Main Cause.
# ---------------------------------------------------------------------
$treeRow = 1
$lightningRow = 1
$lastCol = 14

for ($r = 2; $r -le 40; $r++) {
    $cause = $ws.Cells.Item($r, 3).Value2
    $fval = $ws.Cells.Item($r, 6).Value2
    if ($fval -eq $null) {
        continue
    }
    if ($cause -eq "Tree" -or $cause -eq "Tree/Lightning") {
        $destSheet = $treesSheet
        $destRow = $treeRow
        $treeRow = $treeRow + 1
    } elseif ($cause -eq "Lightning") {
        $destSheet = $lightningSheet
        $destRow = $lightningRow
        $lightningRow = $lightningRow + 1
    } else {
        continue
    }

    for ($c = 1; $c -le $lastCol; $c++) {
        $srcCell = $ws.Cells.Item($r, $c)
        $v = $srcCell.Value2
        if ($v -ne $null) {
            $destSheet.Cells.Item($destRow, $c).Value = $v
        }
    }
}

$treesLastDataRow = $treeRow - 1
$lightningLastDataRow = $lightningRow - 1

# ---------------------------------------------------------------------
# 3. Average summary rows.
# ---------------------------------------------------------------------

# Sheet1: row 43, columns D and E only.
$ws.Range("A43").Value = "Average"
$ws.Range("D43").Formula = "=AVERAGE(D2:D40)"
$ws.Range("E43").Formula = "=AVERAGE(E2:E40)"

# Trees sheet: two rows below the last data row, columns D:L.
$treesAvgRow = $treesLastDataRow + 2
$treesSheet.Cells.Item($treesAvgRow, 1).Value = "Average"
for ($c = 4; $c -le 12; $c++) {
    $colLetter = $treesSheet.Cells.Item(1, $c).Address($false, $false)
    $colLetter = $colLetter -replace '[0-9]', ''
    $cell = $treesSheet.Cells.Item($treesAvgRow, $c)
    $cell.Formula = "=AVERAGE(" + $colLetter + "1:" + $colLetter + $treesLastDataRow + ")"
}

# Lightning sheet: two rows below the last data row, columns D:L.
$lightningAvgRow = $lightningLastDataRow + 2
$lightningSheet.Cells.Item($lightningAvgRow, 1).Value = "Average"
for ($c = 4; $c -le 12; $c++) {
    $colLetter = $lightningSheet.Cells.Item(1, $c).Address($false, $false)
    $colLetter = $colLetter -replace '[0-9]', ''
    $cell = $lightningSheet.Cells.Item($lightningAvgRow, $c)
    $cell.Formula = "=AVERAGE(" + $colLetter + "1:" + $colLetter + $lightningLastDataRow + ")"
}

# ---------------------------------------------------------------------
# 4. View/selection tidy-up on Sheet1 (matches the authored workbook).
# ---------------------------------------------------------------------
$ws.Range("K25").Select()

$win = $excel.ActiveWindow
$win.Left = -31980
$win.Top = -3840
